# Insert a new "groupprivate" column before the existing "grouppublic"
# column (HB), shifting grouppublic/groupwest/group/... and all following
# columns one place to the right (HB->HC, HC->HD, HD->HE, ... , JA->JB).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("HB:HB").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 210).Value = "groupprivate"

# Column numbers (1-based) after the insert:
#   GW=205 group        GX=206 groupsize   GY=207 groupadmit
#   GZ=208 groupwomen   HA=209 groupopen
#   HB=210 groupprivate (new)   HC=211 grouppublic (old HB)
#   HD=212 groupwest (old HC)   HE=213 group (combined label, old HD)
$colGroup       = 205
$colGroupSize   = 206
$colGroupAdmit  = 207
$colGroupWomen  = 208
$colGroupOpen   = 209
$colGroupPriv   = 210
$colGroupPublic = 211
$colGroupWest   = 212
$colGroupLabel  = 213

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $public = $ws.Cells.Item($r, $colGroupPublic).Value()

    if ($public -eq "public") {
        $private = ""
    } else {
        $private = "private"
        $ws.Cells.Item($r, $colGroupPriv).Value = "private"
    }

    $grp      = $ws.Cells.Item($r, $colGroup).Value()
    $size     = $ws.Cells.Item($r, $colGroupSize).Value()
    $admit    = $ws.Cells.Item($r, $colGroupAdmit).Value()
    $women    = $ws.Cells.Item($r, $colGroupWomen).Value()
    $open     = $ws.Cells.Item($r, $colGroupOpen).Value()
    $west     = $ws.Cells.Item($r, $colGroupWest).Value()

    if ([string]::IsNullOrEmpty($size)) { $size = "NA" }
    if ([string]::IsNullOrEmpty($admit)) { $admit = "NA" }

    $label = "$grp $size $admit $women $open $private $public $west"
    $ws.Cells.Item($r, $colGroupLabel).Value = $label
}
